$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("LiveData")

$ws.Range("B2").Value = 7187
$ws.Range("C3").Value = 167072
$ws.Range("C4").Value = 157972
$ws.Range("C8").Value = 65.25
